# Apply the recorded edits to both worksheets:
#  - Sheet1: A2 1 -> 3, A3 2 -> 4, selection moves from E3 to D15
#  - Sheet2: A2 1 -> 3, A3 2 -> 4, selection moves from C9 to A3 (stays the active/tab-selected sheet)

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# --- Sheet1 ---
$ws1.Range("A2").Value = 3
$ws1.Range("A3").Value = 4
$ws1.Range("D15").Select()

# --- Sheet2 (also re-activate it so it remains the selected/active tab) ---
$ws2.Range("A2").Value = 3
$ws2.Range("A3").Value = 4
$ws2.Activate()
$ws2.Range("A3").Select()
